$d = $word.ActiveDocument

# Locate the "Common plan:" paragraph so we can insert the new bullet
# right after it (as the new first item of the existing numbered list).
$found = $d.Content.Find.Execute("Common plan:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$commonPlanPara = $d.Content.Paragraphs(1)
$firstListItem = $commonPlanPara.Next()

# Insert a new paragraph just before the current first list item. It
# inherits that paragraph's formatting (ListParagraph style, numId 7
# numbering), matching the rest of the bulleted list.
$newPara = $firstListItem.Range.InsertParagraphBefore()

# The newly inserted paragraph is now the one right after "Common plan:".
# Set its text accordingly.
$d.Paragraphs(2).Range.Text = "Fix loading (test_deletion1)"
